# Applies the "Making the exercise explicit" commit:
#  1. Adds a new slide (slide 10) "Serving up data on the Interwebs"
#  2. Fixes capitalisation of "program A" -> "Program A" on the
#     "Version Conflicts, System conflicts" slide
#  3. Renames the "lesson7" virtualenv example to "pythonL2" on the
#     "Creating a new virtual environment" slide

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. New slide 10 - "Serving up data on the Interwebs"
# ---------------------------------------------------------------------------
$s = $p.Slides.Add(10, 2)

# --- Title -----------------------------------------------------------------
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Serving up data on the "
[void]$title.InsertAfter("Interwebs")
$title.ParagraphFormat.Alignment = 2

# --- Body content ------------------------------------------------------------
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = "We want to be able to list the restaurants we have data about" + "`r" + "Implement the ListRestaurants method from the previous slide. " + "`r" + "You will need shared.py from a previous week" + [char]0x2019 + "s exercises" + "`r" + "`r"

# Paragraph 1: "We want to be able to list the restaurants we have data about"
# (left as a single run, default formatting)

# Paragraph 2: "Implement the " + "ListRestaurants" (Consolas) + " method from the previous slide. "
$para2 = $body.Paragraphs(2, 1)
$run2code = $para2.Characters(15, 15)
$run2code.Font.Name = "Consolas"

# Paragraph 3: "You will need " + "shared.py" (Consolas) + " " (Consolas) + "from a previous week's exercises"
$para3 = $body.Paragraphs(3, 1)
$run3code = $para3.Characters(15, 9)
$run3code.Font.Name = "Consolas"
$run3space = $para3.Characters(24, 1)
$run3space.Font.Name = "Consolas"

# Paragraph 4: blank, indent level 2 (lvl="1")
$para4 = $body.Paragraphs(4, 1)
$para4.Text = ""
$para4.IndentLevel = 2

# Paragraph 5: blank, indent level 2 (lvl="1"), bullet removed
$para5 = $body.Paragraphs(5, 1)
$para5.Text = ""
$para5.IndentLevel = 2
$para5.ParagraphFormat.Bullet.Type = 0

# --- Decorative logo picture (re-uses the Expedia Code Academy image that
#     already sits in the bottom-right corner of every other slide) --------
$logoPath = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.Shapes.Count -gt 0) {
        foreach ($shp in $candidate.Shapes) {
            if ($shp.Type -eq 13) {
                $logoPath = $shp
                break
            }
        }
    }
    if ($logoPath -ne $null) { break }
}

$x = 9801412 / 12700
$y = 4624575 / 12700
$cx = 1552388 / 12700
$cy = 1552388 / 12700
$pic = $s.Shapes.AddPicture("/tmp/work/extracted/ppt/media/image3.png", $false, $true, $x, $y, $cx, $cy)
$pic.Name = "Picture 4"
$pic.LockAspectRatio = -1

# ---------------------------------------------------------------------------
# 2. "Version Conflicts, System conflicts" slide - capitalise "Program A"
# ---------------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$slide5content = $slide5.Shapes.Item(2).TextFrame.TextRange
$slide5content.Paragraphs(2, 1).Runs(1).Text = "Program A needs version 1 of a module"

# ---------------------------------------------------------------------------
# 3. "Creating a new virtual environment" slide - lesson7 -> pythonL2
# ---------------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$slide6content = $slide6.Shapes.Item(2).TextFrame.TextRange

$slide6content.Paragraphs(2, 1).Runs(2).Text = " pythonL2"
$slide6content.Paragraphs(3, 1).Runs(1).Text = "source pythonL2/bin/activate"
$slide6content.Paragraphs(6, 1).Runs(2).Text = " pythonL2"
$slide6content.Paragraphs(7, 1).Runs(1).Text = "pythonL2\Scripts\activate"
